$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the latest Argent price entry as row 90. The preceding rows
# store the date and price as plain text (not real dates/numbers), so
# force text entry with a leading apostrophe just like typing it in
# Excel, keeping the same representation as rows 74-89.
$ws.Range("A90").Value = "'2025-01-27"
$ws.Range("B90").Value = "'5.83"
